# Fruta / hortaliza, semanal
#
# New weekly price-reporting rows were captured for the "Macroferia Regional
# de Talca - Frutilla" series. Three new records (fecha 44918) are inserted
# right before the existing row 648, pushing every subsequent row down by
# three (old 648..734 -> new 651..737) and growing the used range from
# A1:T734 to A1:T737.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert three blank rows at row 648; this shifts rows 648-734 down to 651-737
# and Excel auto-extends the sheet dimension to A1:T737.
$ws.Rows.Item(648).Insert()
$ws.Rows.Item(648).Insert()
$ws.Rows.Item(648).Insert()

# Columns that are constant across this whole block of "Frutilla" rows.
$constant = @{
    1  = 5                       # A Mercado ID
    2  = "Macroferia Regional de Talca"  # B Mercado
    3  = "Maule"                 # C Región
    5  = 7                       # E Codreg
    6  = "Fruta"                 # F Tipo
    7  = 100101                  # G Producto ID
    8  = "Berries"                # H Producto
    9  = 100112025               # I Categoría ID
    10 = "Frutilla"               # J Categoría
    11 = "Sin especificar"        # K Variedad
    20 = 7                       # T Kg / unidad
}

# Per-row values for the three new records (fecha 44918).
$newRows = @(
    @{ Row = 648; D = 44918; L = "Especial"; M = 50;  N = 7000; O = 7000; P = 7000; Q = "$/bandeja 7 kilos"; R = "Provincia de Melipilla"; S = 1000 },
    @{ Row = 649; D = 44918; L = "Especial"; M = 250; N = 8000; O = 8000; P = 8000; Q = "$/caja 7 kilos";    R = "Región del Maule";        S = 1143 },
    @{ Row = 650; D = 44918; L = "Segunda";  M = 60;  N = 4000; O = 4000; P = 4000; Q = "$/bandeja 7 kilos"; R = "Provincia de Melipilla"; S = 571 }
)

foreach ($rowData in $newRows) {
    $r = $rowData.Row

    foreach ($col in $constant.Keys) {
        $ws.Cells.Item($r, $col).Value = $constant[$col]
    }

    $ws.Cells.Item($r, 4).Value = $rowData.D    # D Fecha
    $ws.Cells.Item($r, 12).Value = $rowData.L   # L Calidad
    $ws.Cells.Item($r, 13).Value = $rowData.M   # M Volumen
    $ws.Cells.Item($r, 14).Value = $rowData.N   # N Precio mínimo
    $ws.Cells.Item($r, 15).Value = $rowData.O   # O Precio máximo
    $ws.Cells.Item($r, 16).Value = $rowData.P   # P Precio promedio ponderado
    $ws.Cells.Item($r, 17).Value = $rowData.Q   # Q Unidad de comercialización
    $ws.Cells.Item($r, 18).Value = $rowData.R   # R Origen
    $ws.Cells.Item($r, 19).Value = $rowData.S   # S Precio $/Kg
}
